# Update the public EPEX Spot prices workbook:
#  - "Prix Spot" sheet: add a new day column AX ("02-aug") with the
#    corresponding hourly prices for rows 2..25 (mirrors the style of AW).
#  - "Gaz" sheet: append a new row 47 for 2025-07-31.
#  - "CO2" sheet: append a new row 47 for 2025-07-31.

$wb = $excel.ActiveWorkbook

### 1) Prix Spot ---------------------------------------------------------
$ws = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the previous day's header cell (AW1) onto the
# new header cell (AX1, column 50), then set its text.
$ws.Range("AW1").Copy($ws.Range("AX1"))
$ws.Range("AX1").Value = "02-aug"

$prixSpotValues = @{
    2  = 97.8
    3  = 77.08
    4  = 65.27
    5  = 49.56
    6  = 43.33
    7  = 44.29
    8  = 51.08
    9  = 49.81
    10 = 37.33
    11 = 23.91
    12 = 16.54
    13 = 5.68
    14 = 2.95
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 8.59
    20 = 34.31
    21 = 65.01000000000001
    22 = 74.53
    23 = 99.64
    24 = 99.64
    25 = 80.5
}

foreach ($row in $prixSpotValues.Keys) {
    $ws.Cells.Item($row, 50).Value = $prixSpotValues[$row]
}

### 2) Gaz ----------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
# Force the date-like string to be stored as plain text (not auto-converted
# to a date serial number), then drop back to the default "Normal" style so
# the cell keeps the same (unstyled) look as the rest of column A.
$wsGaz.Range("A47").NumberFormat = "@"
$wsGaz.Range("A47").Value = "2025-07-31"
$wsGaz.Range("A47").Style = "Normal"
$wsGaz.Range("B47").Value = 34.025

### 3) CO2 ------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A47").NumberFormat = "@"
$wsCo2.Range("A47").Value = "2025-07-31"
$wsCo2.Range("A47").Style = "Normal"
$wsCo2.Range("B47").Value = 71.73999999999999
